$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Repollo / Crespo record was reported, so it
# is inserted at the top of its date-ordered block (row 660). Excel's
# InsertRow shifts every subsequent row (and its formatting) down by one,
# which is exactly what the diff shows happening through row 697.
$ws.Rows.Item(660).Insert()

$ws.Range("A660").Value = 10
$ws.Range("B660").Value = "Vega Modelo de Temuco"
$ws.Range("C660").Value = "La Araucanía"
$ws.Range("D660").Value = 44706
$ws.Range("E660").Value = 9
$ws.Range("F660").Value = 100112006
$ws.Range("G660").Value = "Repollo"
$ws.Range("H660").Value = "Crespo record"
$ws.Range("I660").Value = "Primera"
$ws.Range("J660").Value = 500
$ws.Range("K660").Value = 1400
$ws.Range("L660").Value = 1500
$ws.Range("M660").Value = 1440
$ws.Range("N660").Value = "$/unidad"
$ws.Range("O660").Value = "Región Metropolitana"
$ws.Range("P660").Value = 1440
$ws.Range("Q660").Value = 1
$ws.Range("R660").Value = "Hortaliza"
